$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values E2:T2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.164924666666667
$ws.Range("H2").Value = 3.494774
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.135712
$ws.Range("N2").Value = 0.407136
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1580942563626667
$ws.Range("R2").Value = 1.422848307264
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Delete row 3 entirely
$ws.Rows("3").Delete()

Write-Output "done"
